$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: phone was stored as text ("79174460"); correct it to a real number.
$ws.Range("A6").Value = 79174460

# Row 7: new payment 79174449 (Cash) 2025-08-20T09:46:10
# Column A (phone) must stay text even though it looks numeric, like the
# other phone values in this sheet. Entering a formula that evaluates to a
# text string, then copy/paste-special-values it, converts the cell to a
# plain text value without attaching a new NumberFormat style to it.
$ws.Range("A7").Formula = "=""79174449"""
$ws.Range("A7").Copy()
$ws.Range("A7").PasteSpecial(-4163)

$ws.Range("B7").Value = 70
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 70
$ws.Range("G7").Value = "Cash"
$ws.Range("H7").Value = "2025-08-20T09:46:10"
